$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Module")

# Insert a new column before column E (the "id" column), shifting
# id/name/description one column to the right, then set the header
# of the newly inserted column.
$ws.Range("E1").EntireColumn.Insert()
$ws.Range("E1").Value = "is_metric"
